# g-tail confirmed working. Now to implement the others.
# Adds the "g_union" (and placeholder p_union/q_union/Q_union/single_quote_union)
# lookup columns to the letter table, plus the matching view-state tweaks
# (freeze the first two columns, scroll the header area over, select K67).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (K1:O1) -------------------------------------------------
# Order matters: it drives the shared-string table insertion order.
$ws.Range("K1").Value = "g_union"
$ws.Range("L1").Value = "p_union"
$ws.Range("M1").Value = "q_union"
$ws.Range("N1").Value = "Q_union"
$ws.Range("O1").Value = "single_quote_union"

# --- New column O width (matches column F's "10" character width) ----------
$ws.Columns(15).ColumnWidth = 10

# --- g_union data for rows 2-75 (column K) ----------------------------------
$gUnion = 29,81,29,81,29,82,29,81,83,81,81,84,29,29,29,29,29,29,29,29,83,29,29,29,29,29,85,86,85,86,87,87,85,88,87,84,88,81,89,88,85,86,85,86,85,87,88,88,90,88,88,87,85,83,85,85,81,87,84,87,85,85,29,29,83,83,29,91,92,93,94,29,84,84

for ($i = 0; $i -lt $gUnion.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $gUnion[$i]
}

# --- View state: freeze first two columns, scroll right, select K67 --------
$ws.Range("C1").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 45
$ws.Range("K67").Select()
